# "fixed script for run NF"
# The Scene.xlsx FilePath column (F) referenced scene XML files two
# directories up ("../../NFDataCfg/Ini/Scene/N.xml"); the run-time layout
# only needs one level up, so drop the extra "../" segment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")

$ws.Range("F10").Value = "../NFDataCfg/Ini/Scene/1.xml"
$ws.Range("F11").Value = "../NFDataCfg/Ini/Scene/2.xml"
$ws.Range("F12").Value = "../NFDataCfg/Ini/Scene/3.xml"
$ws.Range("F13").Value = "../NFDataCfg/Ini/Scene/4.xml"
$ws.Range("F14").Value = "../NFDataCfg/Ini/Scene/5.xml"
$ws.Range("F15").Value = "../NFDataCfg/Ini/Scene/6.xml"

# Restore the cursor position left behind in the saved file.
$ws.Activate() | Out-Null
$ws.Cells.Item(23, 6).Select() | Out-Null
